# "updated prelim design report"
#
# Functional change: the "SKR Board" line in the amperage/part-load table
# had its Amperage/unit (column C) corrected from 0.5 to 0.4. The
# dependent formulas (E24 = C24*D24, M24 = B24*C24, and the downstream
# totals in E32/F32) recompute automatically.
#
# The rest of the source diff (fileVersion/rupBuild, absPath UNC, the
# revisionPtr GUID, window placement, and the many small row-height /
# dyDescent tweaks) is just Excel re-rendering the workbook on a
# different machine/display — not an addressable, content-level edit —
# so it is intentionally left alone here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the Amperage/unit for the "SKR Board" row.
$ws.Range("C24").Value = 0.4

# Leave the cursor where the author left it when they saved.
$ws.Range("K30").Select()
